$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 56:57 - shifts existing rows 56:148 down to 58:150
$ws.Rows("56:57").Insert()

# Fill in row 56 (new weekly record - Primera, "$/caja 18 kilos")
$ws.Range("A56").Value = 9
$ws.Range("B56").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C56").Value = "Metropolitana"
$ws.Range("D56").Value = 45012
$ws.Range("E56").Value = 13
$ws.Range("F56").Value = 100114002
$ws.Range("G56").Value = "Camote"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 520
$ws.Range("K56").Value = 17000
$ws.Range("L56").Value = 18000
$ws.Range("M56").Value = 17500
$ws.Range("N56").Value = "$/caja 18 kilos"
$ws.Range("O56").Value = "Perú"
$ws.Range("P56").Value = 972
$ws.Range("Q56").Value = 18
$ws.Range("R56").Value = "Hortaliza"

# Fill in row 57 (new weekly record - Primera, "$/malla 18 kilos")
$ws.Range("A57").Value = 9
$ws.Range("B57").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 45012
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = 100114002
$ws.Range("G57").Value = "Camote"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 790
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 11000
$ws.Range("M57").Value = 10494
$ws.Range("N57").Value = "$/malla 18 kilos"
$ws.Range("O57").Value = "Perú"
$ws.Range("P57").Value = 583
$ws.Range("Q57").Value = 18
$ws.Range("R57").Value = "Hortaliza"
